$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new value, and whether the text must be
# force-typed as Text (so digit-only-looking strings like "193.16" are not
# auto-converted to numbers and lose their original formatting, e.g. "1.00").
$updates = @(
    @('D2', '67.105.04', $false),
    @('E2', '  -1.93%  ', $false),
    @('D3', '3.591.15', $false),
    @('E3', '  -3.02%  ', $false),
    @('E4', '  -0.13%  ', $false),
    @('B5', 'Solana', $false),
    @('C5', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', $false),
    @('D5', '193.16', $true),
    @('E5', '  -1.60%  ', $false),
    @('B6', 'BNB', $false),
    @('C6', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', $false),
    @('D6', '574.92', $true),
    @('E6', '  -6.34%  ', $false),
    @('D7', '3.586.99', $false),
    @('E7', '  -3.02%  ', $false),
    @('E8', '  -2.70%  ', $false),
    @('D9', '0.999', $true),
    @('E9', '  -0.29%  ', $false),
    @('D10', '0.678', $true),
    @('E10', '  -6.34%  ', $false),
    @('E11', '  -5.88%  ', $false),
    @('D12', '55.71', $true),
    @('E12', '  -8.84%  ', $false),
    @('D13', '0.0000269', $true),
    @('E13', '  -5.81%  ', $false),
    @('D14', '9.85', $true),
    @('E14', '  -5.41%  ', $false),
    @('D15', '4.158.83', $false),
    @('E15', '  -3.29%  ', $false),
    @('D16', '3.586.35', $false),
    @('E16', '  -3.30%  ', $false),
    @('D17', '0.126', $true),
    @('E17', '  -1.35%  ', $false),
    @('D18', '18.41', $true),
    @('E18', '  -5.41%  ', $false),
    @('D19', '67.025.88', $false),
    @('E19', '  -1.92%  ', $false),
    @('D20', '12.16', $true),
    @('E20', '  -5.22%  ', $false),
    @('E21', '  -7.18%  ', $false),
    @('D22', '399.71', $true),
    @('E22', '  -2.09%  ', $false),
    @('D23', '4.20', $true),
    @('E23', '  -9.09%  ', $false),
    @('D24', '85.84', $true),
    @('E24', '  -4.48%  ', $false),
    @('D25', '11.46', $true),
    @('E25', '  -0.67%  ', $false),
    @('D26', '2.95', $true),
    @('E26', '  -3.94%  ', $false),
    @('D27', '12.49', $true),
    @('E27', '  -4.60%  ', $false),
    @('E28', '  +0.94%  ', $false),
    @('D29', '3.65', $true),
    @('E29', '  -2.93%  ', $false),
    @('D30', '8.95', $true),
    @('E30', '  -7.03%  ', $false),
    @('D31', '7.67', $true),
    @('E31', '  -0.63%  ', $false),
    @('D32', '31.19', $true),
    @('E32', '  -4.69%  ', $false),
    @('D33', '632.32', $true),
    @('E33', '  -0.42%  ', $false),
    @('D34', '12.18', $true),
    @('E34', '  -3.94%  ', $false),
    @('D35', '0.115', $true),
    @('E35', '  -6.02%  ', $false),
    @('D36', '63.96', $true),
    @('E36', '  -4.99%  ', $false),
    @('D37', '42.12', $true),
    @('E37', '  -12.30%  ', $false),
    @('D38', '0.400', $true),
    @('E38', '  -3.08%  ', $false),
    @('E39', '  +0.09%  ', $false),
    @('D40', '0.0₃0764', $false),
    @('E40', '  -5.68%  ', $false),
    @('D41', '3.181.00', $false),
    @('E41', '  +8.56%  ', $false),
    @('E42', '  -3.75%  ', $false),
    @('D43', '0.999', $true),
    @('E43', '  -0.20%  ', $false),
    @('D44', '2.72', $true),
    @('E44', '  +3.86%  ', $false),
    @('D45', '2.98', $true),
    @('E45', '  -1.86%  ', $false),
    @('D46', '0.0415', $true),
    @('E46', '  -6.70%  ', $false),
    @('B47', 'Stellar', $false),
    @('C47', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', $false),
    @('D47', '0.131', $true),
    @('E47', '  -6.66%  ', $false),
    @('B48', 'ApeXProtocol', $false),
    @('C48', 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex', $false),
    @('D48', '3.06', $true),
    @('E48', '  +0.09%  ', $false),
    @('D49', '140.74', $true),
    @('E49', '  -3.56%  ', $false),
    @('B50', 'dogwifhat', $false),
    @('C50', 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif', $false),
    @('D50', '2.57', $true),
    @('E50', '  -3.23%  ', $false),
    @('B51', 'THORChain', $false),
    @('C51', 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune', $false),
    @('D51', '8.52', $true),
    @('E51', '  -9.26%  ', $false)
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newValue = $u[1]
    $forceText = $u[2]
    $rng = $ws.Range($cellRef)
    if ($forceText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $newValue
}
